# Simplify logging system configuration
# Converts trailing inline-string timestamps in column A to proper numeric
# Excel date-serial values (format "YYYY-MM-DD HH:MM:SS") and appends the
# newest sensor log rows captured on 2025-03-07 to each of the four sheets.
$wb = $excel.ActiveWorkbook

# ---- ROW35-FE-LIFTER ----
$ws = $wb.Worksheets.Item('ROW35-FE-LIFTER')

# Column A: inline-string timestamp -> numeric date serial
$ws.Range("A65").Value = 45721.72894113426
$ws.Range("A65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A66").Value = 45721.72896428241
$ws.Range("A66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A67").Value = 45721.72898766203
$ws.Range("A67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A68").Value = 45722.22908445602
$ws.Range("A68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A69").Value = 45722.22910648148
$ws.Range("A69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A70").Value = 45722.22912973379
$ws.Range("A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A71").Value = 45723.19113143518
$ws.Range("A71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A72").Value = 45723.19115481481
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A73").Value = 45723.19117797454
$ws.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Newly logged rows
$ws.Range("A74").Value = '2025-03-07 16:35:17'
$ws.Range("B74").Value = '0x01,0x90'
$ws.Range("C74").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,'
$ws.Range("D74").Value = '0x01,0x90,'
$ws.Range("E74").Value = '0xd'
$ws.Range("F74").Value = 400
$ws.Range("G74").Value = [double]'5.686312626471138e+23'
$ws.Range("H74").Value = 400
$ws.Range("I74").Value = 13


# ---- ROW35-MID-LIFTER ----
$ws = $wb.Worksheets.Item('ROW35-MID-LIFTER')

# Column A: inline-string timestamp -> numeric date serial
$ws.Range("A68").Value = 45721.72988806713
$ws.Range("A68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A69").Value = 45721.72991133102
$ws.Range("A69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A70").Value = 45721.72993467592
$ws.Range("A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A71").Value = 45722.23020512731
$ws.Range("A71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A72").Value = 45722.23022724537
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A73").Value = 45722.23025050926
$ws.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A74").Value = 45723.19124234954
$ws.Range("A74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A75").Value = 45723.19126548611
$ws.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Newly logged rows
$ws.Range("A76").Value = 45723.19128880787
$ws.Range("A76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B76").Value = '0x01,0x90'
$ws.Range("C76").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,'
$ws.Range("D76").Value = '0x01,0x90,'
$ws.Range("E76").Value = '0xe'
$ws.Range("F76").Value = 400
$ws.Range("G76").Value = [double]'5.686312626471138e+23'
$ws.Range("H76").Value = 400
$ws.Range("I76").Value = 14

$ws.Range("A77").Value = '2025-03-07 16:35:27'
$ws.Range("B77").Value = '0x01,0x90'
$ws.Range("C77").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,'
$ws.Range("D77").Value = '0x01,0x90,'
$ws.Range("E77").Value = '0xe'
$ws.Range("F77").Value = 400
$ws.Range("G77").Value = [double]'5.686312626471138e+23'
$ws.Range("H77").Value = 400
$ws.Range("I77").Value = 14


# ---- ROW02-FE-LIFTER ----
$ws = $wb.Worksheets.Item('ROW02-FE-LIFTER')

# Column A: inline-string timestamp -> numeric date serial
$ws.Range("A65").Value = 45721.72979079861
$ws.Range("A65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A66").Value = 45721.72981402778
$ws.Range("A66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A67").Value = 45721.7298375463
$ws.Range("A67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A68").Value = 45722.23010768519
$ws.Range("A68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A69").Value = 45722.23012966435
$ws.Range("A69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A70").Value = 45722.23015302084
$ws.Range("A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A71").Value = 45723.19132489583
$ws.Range("A71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A72").Value = 45723.19134847222
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A73").Value = 45723.19137149306
$ws.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Newly logged rows
$ws.Range("A74").Value = '2025-03-07 16:35:34'
$ws.Range("B74").Value = '0x01,0x90'
$ws.Range("C74").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,'
$ws.Range("D74").Value = '0x01,0x90,'
$ws.Range("E74").Value = '0xff'
$ws.Range("F74").Value = 400
$ws.Range("G74").Value = [double]'5.686312626471138e+23'
$ws.Range("H74").Value = 400
$ws.Range("I74").Value = 255


# ---- ROW02-MID-LIFTER ----
$ws = $wb.Worksheets.Item('ROW02-MID-LIFTER')

# Column A: inline-string timestamp -> numeric date serial
$ws.Range("A65").Value = 45721.72820228009
$ws.Range("A65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A66").Value = 45721.7282258449
$ws.Range("A66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A67").Value = 45721.72824888889
$ws.Range("A67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A68").Value = 45722.22834673611
$ws.Range("A68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A69").Value = 45722.22836789352
$ws.Range("A69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A70").Value = 45722.22839168982
$ws.Range("A70").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A71").Value = 45722.72848770834
$ws.Range("A71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A72").Value = 45722.72850997685
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A73").Value = 45722.72853335648
$ws.Range("A73").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Newly logged rows
$ws.Range("A74").Value = 45723.22863142361
$ws.Range("A74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B74").Value = '0x01,0x90'
$ws.Range("C74").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,'
$ws.Range("D74").Value = '0x01,0x90,'
$ws.Range("E74").Value = '0x3'
$ws.Range("F74").Value = 400
$ws.Range("G74").Value = [double]'5.686312626471138e+23'
$ws.Range("H74").Value = 400
$ws.Range("I74").Value = 3

$ws.Range("A75").Value = 45723.22865329861
$ws.Range("A75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B75").Value = '0x01,0x90'
$ws.Range("C75").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,'
$ws.Range("D75").Value = '0x01,0x90,'
$ws.Range("E75").Value = '0x3'
$ws.Range("F75").Value = 400
$ws.Range("G75").Value = [double]'5.686312626471138e+23'
$ws.Range("H75").Value = 400
$ws.Range("I75").Value = 3

$ws.Range("A76").Value = 45723.22867643519
$ws.Range("A76").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B76").Value = '0x01,0x90'
$ws.Range("C76").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,'
$ws.Range("D76").Value = '0x01,0x90,'
$ws.Range("E76").Value = '0x3'
$ws.Range("F76").Value = 400
$ws.Range("G76").Value = [double]'5.686312626471138e+23'
$ws.Range("H76").Value = 400
$ws.Range("I76").Value = 3

$ws.Range("A77").Value = '2025-03-07 17:29:17'
$ws.Range("B77").Value = '0x01,0x90'
$ws.Range("C77").Value = '0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,'
$ws.Range("D77").Value = '0x01,0x90,'
$ws.Range("E77").Value = '0x3'
$ws.Range("F77").Value = 400
$ws.Range("G77").Value = [double]'5.686312626471138e+23'
$ws.Range("H77").Value = 400
$ws.Range("I77").Value = 3

